$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("B1").Value = 0.2028864896758006
$ws.Range("C1").Value = 0.002539495620456482
$ws.Range("D1").Value = -1.379316785363834
$ws.Range("E1").Value = 0.188939888900502
$ws.Range("F1").Value = 1.570796384046431

# Row 2
$ws.Range("B2").Value = 0.2484977296296271
$ws.Range("C2").Value = 0.002377709119405816
$ws.Range("D2").Value = -1.380518500463224
$ws.Range("E2").Value = 0.1878999627744357
$ws.Range("F2").Value = 1.570796384423863

# Row 3
$ws.Range("B3").Value = 0.4528663241188721
$ws.Range("C3").Value = 0.001652798332930176
$ws.Range("D3").Value = -1.385902980825687
$ws.Range("E3").Value = 0.1832404042758736
$ws.Range("F3").Value = 1.570796386115009

# Row 4
$ws.Range("B4").Value = 0.7403935328885739
$ws.Range("C4").Value = 0.0006329176826605036
$ws.Range("D4").Value = -1.393478433468806
$ws.Range("E4").Value = 0.1766848478679089
$ws.Range("F4").Value = 1.570796388494291

# Row 5
$ws.Range("B5").Value = 0.9447621273778194
$ws.Range("C5").Value = -0.00009199310381513776
$ws.Range("D5").Value = -1.398862913831269
$ws.Range("E5").Value = 0.1720252893693467
$ws.Range("F5").Value = 1.570796390185437

# Row 6
$ws.Range("B6").Value = 0.9903733673316449
$ws.Range("C6").Value = -0.0002537796048658041
$ws.Range("D6").Value = -1.40006462893066
$ws.Range("E6").Value = 0.1709853632432805
$ws.Range("F6").Value = 1.570796390562869
